$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.178.67"
$ws.Range("D2").Style = 'Normal'
$ws.Range("E2").Value = '  -2.85%  '
$ws.Range("D3").Value = "'1.712.68"
$ws.Range("D3").Style = 'Normal'
$ws.Range("E3").Value = '  -3.38%  '
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").Value = "'308.40"
$ws.Range("D5").Style = 'Normal'
$ws.Range("E5").Value = '  -6.09%  '
$ws.Range("E6").Value = '  +0.26%  '
$ws.Range("D7").Value = "'0.4744"
$ws.Range("D7").Style = 'Normal'
$ws.Range("E7").Value = '  +5.74%  '
$ws.Range("D8").Value = "'0.3442"
$ws.Range("D8").Style = 'Normal'
$ws.Range("E8").Value = '  -3.28%  '
$ws.Range("D9").Value = "'42.02"
$ws.Range("D9").Style = 'Normal'
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("D10").Value = "'0.07287"
$ws.Range("D10").Style = 'Normal'
$ws.Range("E10").Value = '  -1.93%  '
$ws.Range("D11").Value = "'1.045"
$ws.Range("D11").Style = 'Normal'
$ws.Range("E11").Value = '  -5.55%  '
$ws.Range("D13").Value = "'19.87"
$ws.Range("D13").Style = 'Normal'
$ws.Range("E13").Value = '  -5.33%  '
$ws.Range("D14").Value = "'5.861"
$ws.Range("D14").Style = 'Normal'
$ws.Range("E14").Value = '  -2.98%  '
$ws.Range("D15").Value = "'1.709.65"
$ws.Range("D15").Style = 'Normal'
$ws.Range("E15").Value = '  -3.22%  '
$ws.Range("D16").Value = "'6.850"
$ws.Range("D16").Style = 'Normal'
$ws.Range("E16").Value = '  -5.49%  '
$ws.Range("D17").Value = "'88.79"
$ws.Range("D17").Style = 'Normal'
$ws.Range("E17").Value = '  -5.08%  '
$ws.Range("D18").Value = "'0.00001041"
$ws.Range("D18").Style = 'Normal'
$ws.Range("E18").Value = '  -2.11%  '
$ws.Range("D19").Value = "'0.06362"
$ws.Range("D19").Style = 'Normal'
$ws.Range("E19").Value = '  -1.10%  '
$ws.Range("E22").Value = '  -3.14%  '
$ws.Range("D23").Value = "'27.218.77"
$ws.Range("D23").Style = 'Normal'
$ws.Range("E23").Value = '  -2.68%  '
$ws.Range("D24").Value = "'10.79"
$ws.Range("D24").Style = 'Normal'
$ws.Range("E24").Value = '  -4.44%  '
$ws.Range("E25").Value = '  -0.75%  '
$ws.Range("D26").Value = "'153.08"
$ws.Range("D26").Style = 'Normal'
$ws.Range("E26").Value = '  -5.24%  '
$ws.Range("D27").Value = "'19.81"
$ws.Range("D27").Style = 'Normal'
$ws.Range("E27").Value = '  -2.67%  '
$ws.Range("D28").Value = "'1.905.69"
$ws.Range("D28").Style = 'Normal'
$ws.Range("E28").Value = '  -3.18%  '
$ws.Range("E29").Value = '  -3.41%  '
$ws.Range("E30").Value = '  -3.61%  '
$ws.Range("D31").Value = "'1.015"
$ws.Range("D31").Style = 'Normal'
$ws.Range("E31").Value = '  -8.60%  '
$ws.Range("D32").Value = "'0.09224"
$ws.Range("D32").Style = 'Normal'
$ws.Range("E32").Value = '  +0.42%  '
$ws.Range("D33").Value = "'3.590"
$ws.Range("D33").Style = 'Normal'
$ws.Range("E33").Value = '  -2.64%  '
$ws.Range("D34").Value = "'5.299"
$ws.Range("D34").Style = 'Normal'
$ws.Range("E34").Value = '  -6.95%  '
$ws.Range("D35").Value = "'0.02194"
$ws.Range("D35").Style = 'Normal'
$ws.Range("E35").Value = '  -4.30%  '
$ws.Range("D36").Value = "'0.05889"
$ws.Range("D36").Style = 'Normal'
$ws.Range("E36").Value = '  -5.30%  '
$ws.Range("E37").Value = '  -6.85%  '
$ws.Range("D38").Value = "'0.2008"
$ws.Range("D38").Style = 'Normal'
$ws.Range("E38").Value = '  -4.84%  '
$ws.Range("D39").Value = "'4.741"
$ws.Range("D39").Style = 'Normal'
$ws.Range("E39").Value = '  -4.73%  '
$ws.Range("E40").Value = '  +0.36%  '
$ws.Range("D41").Value = "'1.410"
$ws.Range("D41").Style = 'Normal'
$ws.Range("E41").Value = '  +1.02%  '
$ws.Range("D42").Value = "'0.5927"
$ws.Range("D42").Style = 'Normal'
$ws.Range("E42").Value = '  -6.19%  '
$ws.Range("D43").Value = "'1.113"
$ws.Range("D43").Style = 'Normal'
$ws.Range("E43").Value = '  -6.16%  '
$ws.Range("D44").Value = "'7.472"
$ws.Range("D44").Style = 'Normal'
$ws.Range("E44").Value = '  -5.29%  '
$ws.Range("D45").Value = "'12.69"
$ws.Range("D45").Style = 'Normal'
$ws.Range("E45").Value = '  -4.80%  '
$ws.Range("E46").Value = '  -4.92%  '
$ws.Range("D47").Value = "'0.5622"
$ws.Range("D47").Style = 'Normal'
$ws.Range("E47").Value = '  -4.51%  '
$ws.Range("D48").Value = "'118.40"
$ws.Range("D48").Style = 'Normal'
$ws.Range("E48").Value = '  -3.32%  '
$ws.Range("D49").Value = "'1.844"
$ws.Range("D49").Style = 'Normal'
$ws.Range("E49").Value = '  -5.99%  '
$ws.Range("D50").Value = "'0.06630"
$ws.Range("D50").Style = 'Normal'
$ws.Range("E50").Value = '  -3.70%  '
$ws.Range("E51").Value = '  -4.68%  '
